$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Update the short "Sociodemographic analysis ..." paragraph text.
#    This paragraph is a single run / single <w:t> with no special
#    whitespace, so a plain Range.Text assignment reproduces the
#    target markup exactly.
# ------------------------------------------------------------------
$newIntro = "Sociodemographic analysis was integrated by including Gender and SES (Socioeconomic Status) as conditions. Gender was dummy-coded (1=Female, 0=Male) and SES was dichotomized (1=High SES [Level 3+], 0=Low SES). The sufficiency analysis seeks the minimal combination of empathy dimensions and sociodemographic conditions leading to high empathy."

foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Sociodemographic analysis was integrated*") {
        $p.Range.Text = $newIntro
        break
    }
}

# ------------------------------------------------------------------
# 2) Rebuild the "Table 4" fsQCA results block. This paragraph mixes
#    <w:t> runs that keep default whitespace handling with others
#    that require xml:space="preserve" (leading/trailing spaces),
#    separated by <w:br/> line breaks. Plain Range.Text assignment
#    cannot reproduce that mixed preserve/non-preserve pattern, so we
#    inject the exact OOXML for the run via Range.InsertXML, which
#    replaces the whole paragraph's contents with precisely what we
#    specify (keeping the same Courier New / sz=18 run formatting the
#    paragraph already had).
# ------------------------------------------------------------------
$tableXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:r>' +
    '<w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New"/><w:sz w:val="18"/></w:rPr>' +
    '<w:t>(With Remainders) ---</w:t>' +
    '<w:br/>' +
    '<w:br/>' +
    '<w:t>M1: fs_f*pt_f + fs_f*ec_f + fs_f*pd_f + fs_f*gen_f + pt_f*ec_f + pt_f*pd_f +</w:t>' +
    '<w:br/>' +
    '<w:t xml:space="preserve">    pt_f*gen_f + ec_f*pd_f + ec_f*ses_f + pd_f*gen_f*~ses_f -&gt; iri_total_f</w:t>' +
    '<w:br/>' +
    '<w:br/>' +
    '<w:t xml:space="preserve">                       inclS   PRI   covS   covU  </w:t>' +
    '<w:br/>' +
    '<w:t xml:space="preserve">------------------------------------------------- </w:t>' +
    '<w:br/>' +
    '<w:t xml:space="preserve"> 1          fs_f*pt_f  0.935  0.869  0.677  0.020 </w:t>' +
    '<w:br/>' +
    '<w:t xml:space="preserve"> 2          fs_f*ec_f  0.961  0.923  0.731  0.006 </w:t>' +
    '<w:br/>' +
    '<w:t xml:space="preserve"> 3          fs_f*pd_f  0.918  0.835  0.661  0.018 </w:t>' +
    '<w:br/>' +
    '<w:t xml:space="preserve"> 4         fs_f*gen_f  0.896  0.833  0.455  0.004 </w:t>' +
    '<w:br/>' +
    '<w:t xml:space="preserve"> 5          pt_f*ec_f  0.918  0.839  0.719  0.014 </w:t>' +
    '<w:br/>' +
    '<w:t xml:space="preserve"> 6          pt_f*pd_f  0.940  0.869  0.598  0.005 </w:t>' +
    '<w:br/>' +
    '<w:t xml:space="preserve"> 7         pt_f*gen_f  0.860  0.775  0.424  0.002 </w:t>' +
    '<w:br/>' +
    '<w:t xml:space="preserve"> 8          ec_f*pd_f  0.944  0.886  0.676  0.004 </w:t>' +
    '<w:br/>' +
    '<w:t xml:space="preserve"> 9         ec_f*ses_f  0.871  0.781  0.239  0.001 </w:t>' +
    '<w:br/>' +
    '<w:t xml:space="preserve">10  pd_f*gen_f*~ses_f  0.844  0.747  0.318  0.002 </w:t>' +
    '<w:br/>' +
    '<w:t xml:space="preserve">------------------------------------------------- </w:t>' +
    '<w:br/>' +
    '<w:t xml:space="preserve">                   M1  0.773  0.634  0.984</w:t>' +
  '</w:r>' +
'</w:p>'

foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*With Remainders*") {
        $p.Range.InsertXML($tableXml) | Out-Null
        break
    }
}
